$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'capri mens joggers',
    'easton baseball pants mens',
    'mueller knee pads basketball',
    'bmx knee pads youth',
    'basketball tights kids',
    'tortoise knee pads',
    'bendable knee pads',
    'dakine knee pads',
    '187 knee pads',
    'mma knee pads',
    'scrubs men pants',
    'copper compression pants',
    'elbow pads knee pads',
    'knee pad for scooter',
    'hamstring compression pants',
    'compression tights men basketball',
    'white leggings for men',
    'red leggings men',
    'purple leggings men',
    'yellow leggings men',
    'athletic leggings for women',
    'athletic leggings kids',
    'basketball leggings for women',
    'basketball leggings kids',
    'bdu pants with knee pads',
    'nike leggings for men',
    'compression pants basketball',
    'compression pants for women',
    'compression pants kids',
    'compression pants knee',
    'compression pants men under armour',
    'compression pants set',
    'compression pants tesla',
    'camo pants with knee pads',
    'kids basketball pads',
    'gold leggings men',
    'kickboxing knee pads',
    'men gym pants',
    'men nike compression pants',
    'swim leggings for men',
    'fleece leggings men',
    'sliding shorts with knee pads',
    'mens basketball joggers',
    'mcdavid basketball knee',
    '3xl knee pads',
    '511 knee pads',
    'caving knee pads',
    'enduro knee pads',
    'bodyprox knee pads',
    'blackhawk knee pads',
    'snickers knee pads',
    'bball knee pads',
    '661 knee pads',
    'basket knee pads',
    'armadillo knee pads',
    'swim pants men',
    'athletic capris',
    'bunheads knee pads',
    '6xl compression pants',
    'kali knee pads',
    'arcteryx knee pads',
    'spelunking knee pads',
    'bcg compression pants',
    'bcg knee pads',
    'tesla mens leggings',
    'awp knee pads',
    'leatt knee pads',
    'fr knee pads',
    'eurotard knee pads',
    'asics leggings men',
    'alleson baseball pants youth',
    'frozen knee pads',
    'training tights men',
    'training leggings for men',
    'nba knee pads',
    'prayer knee pads',
    'neoprene knee pad',
    's1 knee pads',
    'white basketball pants',
    'lululemon compression pants',
    'xlarge knee pads',
    'ballet tights for men',
    'nike compression pants for men',
    'mens compression pants adidas',
    'mens compression pants champion',
    'baseball pants mens with piping',
    'exercise compression leggings',
    'teen knee pads',
    'short tights for men',
    'runner compression pants',
    'swimming tights men',
    'half tights men',
    'fishing knee pads',
    'marvel knee pads',
    'bauer compression pants',
    'police knee pads',
    'spiderman tights men',
    'wwe knee pads',
    'petite compression leggings',
    '3x compression leggings'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
